$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q4" right after "总计", using a copy of
#    "2022-Q3" as a template so number formats / styles / column widths are
#    inherited, then overwrite its values with the new quarter's data.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsQ3Old = $wb.Worksheets.Item(2)
$wsQ3Old.Copy($null, $wsTotal)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Columns B..G hold text-formatted values (fund code/name/size/position/
# ratio/market value) even though several look numeric - force text storage
# to match the source data's inline-string cells.
$wsQ4.Range("B2:G5").NumberFormat = "@"

$wsQ4.Range("A2").Value2 = 0
$wsQ4.Range("B2").Value2 = "160416"
$wsQ4.Range("C2").Value2 = "华安标普全球石油指数（QDII-LOF）A"
$wsQ4.Range("D2").Value2 = "2.81"
$wsQ4.Range("E2").Value2 = "93.63"
$wsQ4.Range("F2").Value2 = "5.58"
$wsQ4.Range("G2").Value2 = "0.1568"
$wsQ4.Range("H2").Value2 = 3

$wsQ4.Range("A3").Value2 = 1
$wsQ4.Range("B3").Value2 = "014982"
$wsQ4.Range("C3").Value2 = "华安标普全球石油指数（QDII-LOF）C"
$wsQ4.Range("D3").Value2 = "0.36"
$wsQ4.Range("E3").Value2 = "93.63"
$wsQ4.Range("F3").Value2 = "5.58"
$wsQ4.Range("G3").Value2 = "0.0201"
$wsQ4.Range("H3").Value2 = 3

$wsQ4.Range("A4").Value2 = 2
$wsQ4.Range("B4").Value2 = "010343"
$wsQ4.Range("C4").Value2 = "华宝英国富时100指数A"
$wsQ4.Range("D4").Value2 = "0.14"
$wsQ4.Range("E4").Value2 = "94.75"
$wsQ4.Range("F4").Value2 = "8.22"
$wsQ4.Range("G4").Value2 = "0.0115"
$wsQ4.Range("H4").Value2 = 2

$wsQ4.Range("A5").Value2 = 3
$wsQ4.Range("B5").Value2 = "010344"
$wsQ4.Range("C5").Value2 = "华宝英国富时100指数C"
$wsQ4.Range("D5").Value2 = "0.08"
$wsQ4.Range("E5").Value2 = "94.75"
$wsQ4.Range("F5").Value2 = "8.22"
$wsQ4.Range("G5").Value2 = "0.0066"
$wsQ4.Range("H5").Value2 = 2

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary sheet: add the new 2022-Q4 row at the top of
#    the data and shift the existing quarters down by one row.
# ---------------------------------------------------------------------------
$wsTotal.Range("A5").Value2 = 3
$wsTotal.Range("B5").Value2 = "2022-Q1"
$wsTotal.Range("C5").Value2 = 3
$wsTotal.Range("D5").Value2 = 0.13

$wsTotal.Range("A4").Value2 = 2
$wsTotal.Range("B4").Value2 = "2022-Q2"
$wsTotal.Range("C4").Value2 = 4
$wsTotal.Range("D4").Value2 = 0.23

$wsTotal.Range("A3").Value2 = 1
$wsTotal.Range("B3").Value2 = "2022-Q3"
$wsTotal.Range("C3").Value2 = 4
$wsTotal.Range("D3").Value2 = 0.17

$wsTotal.Range("A2").Value2 = 0
$wsTotal.Range("B2").Value2 = "2022-Q4"
$wsTotal.Range("C2").Value2 = 4
$wsTotal.Range("D2").Value2 = 0.2

# ---------------------------------------------------------------------------
# 3. Keep "2022-Q1" (the original last tab) as the active/selected sheet, as
#    it was in the source workbook before the new sheet was inserted.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
